# Applies the edits described by the commit:
# "Capitulo dos revisado hasta versiculo 10"
# (Chapter two reviewed up to verse 10)
#
# Summary of changes:
#  - Trim trailing spaces from three existing notes.
#  - Personas!Lugar sheet ("Lugar"): add a "Judá; tierra de Judá;" variation
#    for Judea (pla2), and populate the pla39 row (oriente / Pais / comment).
#  - Organizaciones: add "judíos" as a name for org13.
#  - Momento: add "días del Rey Herodes" as a name for tim8.
#  - Switch the active/selected sheet from "Personas" to "Lugar", and update
#    the remembered cell selection on a few sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Text clean-up (trailing spaces removed)
# ---------------------------------------------------------------------
$wsPersonas = $wb.Worksheets.Item("Personas")
$wsPersonas.Range("C12").Value = "María; madre;"
$wsPersonas.Range("C19").Value = "Juan;"

$wsMomento = $wb.Worksheets.Item("Momento")
$wsMomento.Range("B3").Value = "Día del juicio; fin del mundo;"

# ---------------------------------------------------------------------
# 2. New data entered on the "Lugar" sheet
# ---------------------------------------------------------------------
$wsLugar = $wb.Worksheets.Item("Lugar")

# pla2 / Judea row gains a "variaciones" entry
$wsLugar.Range("C3").Value = "Judá; tierra de Judá; "

# pla39 row gains nombre / variaciones / Comentario
$wsLugar.Range("B40").Value = "oriente"
$wsLugar.Range("C40").Value = "País"
$wsLugar.Range("D40").Value = "Lugar del que proceden los magos que alaban a Jesús"

# ---------------------------------------------------------------------
# 3. New data entered on the "Organizaciones" sheet
# ---------------------------------------------------------------------
$wsOrg = $wb.Worksheets.Item("Organizaciones")
$wsOrg.Range("B14").Value = "judíos"

# ---------------------------------------------------------------------
# 4. New data entered on the "Momento" sheet
# ---------------------------------------------------------------------
$wsMomento.Range("B9").Value = "días del Rey Herodes"

# ---------------------------------------------------------------------
# 5. Update selections to match where the reviewer left off
# ---------------------------------------------------------------------
$wsOrg.Range("A13").Select()
$wsMomento.Range("A9").Select()

# Make "Lugar" the active sheet (drives workbook.xml's activeTab and this
# sheet's tabSelected flag), then leave the selection on D40.
$wsLugar.Activate()
$wsLugar.Range("D40").Select()
